$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-186 all get updated to the new date serial value 45202
$ws.Range("C2:C186").Value = 45202
